$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the shared-string text value "CA1402" in A2 with the numeric value 400000002
$ws.Range("A2").Value = 400000002
